# Update row 8 (year 2025) metrics in the BIBI annual recurrence metrics sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1305
$ws.Range("E8").Value = 1098
$ws.Range("G8").Value = 84.13793103448276
$ws.Range("H8").Value = 15.86206896551724
